$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 39
$ws.Cells.Item(39,1).Value = '15/08/2018'
$ws.Cells.Item(39,2).Value = 'B'
$ws.Cells.Item(39,3).Value = 80001841
$ws.Cells.Item(39,5).Value = 5
$ws.Cells.Item(39,6).Value = 'електрична поломка'
$ws.Cells.Item(39,7).Value = 'тест'
# Row 40
$ws.Cells.Item(40,1).Value = '02/01/2019'
$ws.Cells.Item(40,2).Value = '*'
$ws.Cells.Item(40,3).Value = 80001841
$ws.Cells.Item(40,4).Value = '057'
$ws.Cells.Item(40,5).Value = 10
$ws.Cells.Item(40,6).Value = 'scaner'
# Row 41
$ws.Cells.Item(41,1).Value = '02/01/2019'
$ws.Cells.Item(41,2).Value = '*'
$ws.Cells.Item(41,3).Value = 80001841
$ws.Cells.Item(41,5).Value = 10
$ws.Cells.Item(41,6).Value = 'електрична поломка'
# Row 42
$ws.Cells.Item(42,1).Value = '02/01/2019'
$ws.Cells.Item(42,2).Value = '*'
$ws.Cells.Item(42,3).Value = 80001841
$ws.Cells.Item(42,5).Value = 10
$ws.Cells.Item(42,6).Value = 'електрична поломка'
# Row 43
$ws.Cells.Item(43,1).Value = '02/01/2019'
$ws.Cells.Item(43,2).Value = '*'
$ws.Cells.Item(43,3).Value = 80001841
$ws.Cells.Item(43,5).Value = 10
$ws.Cells.Item(43,6).Value = 'заміна запчастин'
# Row 44
$ws.Cells.Item(44,1).Value = '02/01/2019'
$ws.Cells.Item(44,2).Value = '*'
$ws.Cells.Item(44,3).Value = 80001841
$ws.Cells.Item(44,4).Value = '057'
$ws.Cells.Item(44,5).Value = 10
$ws.Cells.Item(44,6).Value = 'електрична поломка'
# Row 45
$ws.Cells.Item(45,1).Value = '02/01/2019'
$ws.Cells.Item(45,2).Value = '*'
$ws.Cells.Item(45,3).Value = 80001841
$ws.Cells.Item(45,4).Value = '057'
$ws.Cells.Item(45,5).Value = 11
$ws.Cells.Item(45,6).Value = 'електрична поломка'
# Row 46
$ws.Cells.Item(46,1).Value = '02/01/2019'
$ws.Cells.Item(46,2).Value = '*'
$ws.Cells.Item(46,3).Value = 80001841
$ws.Cells.Item(46,5).Value = 11
$ws.Cells.Item(46,6).Value = 'механічне налаштування'
# Row 47
$ws.Cells.Item(47,1).Value = '02/01/2019'
$ws.Cells.Item(47,2).Value = '*'
$ws.Cells.Item(47,3).Value = 80001841
$ws.Cells.Item(47,4).Value = '0571'
$ws.Cells.Item(47,5).Value = 11
$ws.Cells.Item(47,6).Value = 'механічне налаштування'
# Row 48
$ws.Cells.Item(48,1).Value = '02/01/2019'
$ws.Cells.Item(48,2).Value = '*'
$ws.Cells.Item(48,3).Value = 80001841
$ws.Cells.Item(48,4).Value = '0571'
$ws.Cells.Item(48,5).Value = 11
$ws.Cells.Item(48,6).Value = 'налаштування втулочного модуля'
# Row 49
$ws.Cells.Item(49,1).Value = '02/01/2019'
$ws.Cells.Item(49,2).Value = '*'
$ws.Cells.Item(49,3).Value = 80001841
$ws.Cells.Item(49,5).Value = 10
$ws.Cells.Item(49,6).Value = 'налаштування симетричності розрізу'
# Row 50
$ws.Cells.Item(50,1).Value = '02/01/2019'
$ws.Cells.Item(50,2).Value = '*'
$ws.Cells.Item(50,3).Value = 80001841
$ws.Cells.Item(50,5).Value = 10
$ws.Cells.Item(50,6).Value = 'заміна запчастин'
# Row 51
$ws.Cells.Item(51,1).Value = '02/01/2019'
$ws.Cells.Item(51,2).Value = '*'
$ws.Cells.Item(51,3).Value = 80001841
$ws.Cells.Item(51,5).Value = 10
$ws.Cells.Item(51,6).Value = 'ПЗ'
# Row 52
$ws.Cells.Item(52,1).Value = '02/01/2019'
$ws.Cells.Item(52,2).Value = '*'
$ws.Cells.Item(52,3).Value = 80001841
$ws.Cells.Item(52,5).Value = 10
$ws.Cells.Item(52,6).Value = 'інший тип простою'
# Row 53
$ws.Cells.Item(53,1).Value = '02/01/2019'
$ws.Cells.Item(53,2).Value = '*'
$ws.Cells.Item(53,3).Value = 80001841
$ws.Cells.Item(53,5).Value = 10
$ws.Cells.Item(53,6).Value = 'ТО аплікатора'
# Row 54
$ws.Cells.Item(54,1).Value = '02/01/2019'
$ws.Cells.Item(54,2).Value = '*'
$ws.Cells.Item(54,3).Value = 80001841
$ws.Cells.Item(54,4).Value = '0571'
$ws.Cells.Item(54,5).Value = 10
$ws.Cells.Item(54,6).Value = 'заміна запчастин'
# Row 55
$ws.Cells.Item(55,1).Value = '02/01/2019'
$ws.Cells.Item(55,2).Value = '*'
$ws.Cells.Item(55,3).Value = 80001841
$ws.Cells.Item(55,4).Value = '0571'
$ws.Cells.Item(55,5).Value = 10
$ws.Cells.Item(55,6).Value = 'механічне налаштування'
# Row 56
$ws.Cells.Item(56,1).Value = '02/01/2019'
$ws.Cells.Item(56,2).Value = '*'
$ws.Cells.Item(56,3).Value = 80001841
$ws.Cells.Item(56,4).Value = '0571'
$ws.Cells.Item(56,5).Value = 10
$ws.Cells.Item(56,6).Value = 'механічне налаштування'
# Row 57
$ws.Cells.Item(57,1).Value = '02/01/2019'
$ws.Cells.Item(57,2).Value = '*'
$ws.Cells.Item(57,3).Value = 80001841
$ws.Cells.Item(57,4).Value = '0571'
$ws.Cells.Item(57,5).Value = 10
$ws.Cells.Item(57,6).Value = 'механічне налаштування'
# Row 58
$ws.Cells.Item(58,1).Value = '02/01/2019'
$ws.Cells.Item(58,2).Value = '*'
$ws.Cells.Item(58,3).Value = 80001841
$ws.Cells.Item(58,5).Value = 10
$ws.Cells.Item(58,6).Value = 'налаштування принтера'
# Row 59
$ws.Cells.Item(59,1).Value = '02/01/2019'
$ws.Cells.Item(59,2).Value = '*'
$ws.Cells.Item(59,3).Value = 80001841
$ws.Cells.Item(59,4).Value = '0571'
$ws.Cells.Item(59,5).Value = 10
$ws.Cells.Item(59,6).Value = 'заміна запчастин'
# Row 60
$ws.Cells.Item(60,1).Value = '02/01/2019'
$ws.Cells.Item(60,2).Value = '*'
$ws.Cells.Item(60,3).Value = 80001841
$ws.Cells.Item(60,5).Value = 10
$ws.Cells.Item(60,6).Value = 'механічна поломка'
# Row 61
$ws.Cells.Item(61,1).Value = '**'
